$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append an additional remark to Todd Robinson's existing P1 comment (row 27)
$ws.Range("B27").Value = " no data process discussion (-1%), late video (-1)"

# Add new P1 comments for students that did not have one yet
$ws.Range("B3").Value = "Late image (-0.5), not telling data process (-1), late video (-1), late readme update(-1),"
$ws.Range("B16").Value = "No data description (-1), failed to show image in readme (-0.5)"
$ws.Range("B17").Value = "No video (-3), No findings (-1)"
$ws.Range("B40").Value = "Late video (-1),late update readme(-1), No Finding (-1)"

# Update the view's selection to match where the editor left off
$ws.Range("B19").Select()
